$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.1467555
$ws.Range("H2").Value = 20.293511
$ws.Range("I2").Value = 0.1516003594919049
$ws.Range("J2").Value = 0.1102643619993968
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.111657
$ws.Range("N2").Value = 0.223314
$ws.Range("O2").Value = 0.0004795212148781231
$ws.Range("P2").Value = 0.0003199968420168896
$ws.Range("Q2").Value = 1.1329562788635
$ws.Range("R2").Value = 4.531825115454001
$ws.Range("S2").Value = 0.00007269558855951843
$ws.Range("T2").Value = 0.00003528424762681409

$ws.Range("G3").Value = 10.1467555
$ws.Range("H3").Value = 20.293511
$ws.Range("I3").Value = 0.1516003594919049
$ws.Range("J3").Value = 0.1102643619993968
$ws.Range("O3").Value = 0.001493600267587414
$ws.Range("P3").Value = 0.001495076820485504
$ws.Range("Q3").Value = 3.528902890574833
$ws.Range("R3").Value = 21.173417343449
$ws.Range("S3").Value = 0.0002264303375034573
$ws.Range("T3").Value = 0.0001648536917509207

$ws.Range("G4").Value = 10.1467555
$ws.Range("H4").Value = 20.293511
$ws.Range("I4").Value = 0.1516003594919049
$ws.Range("J4").Value = 0.1102643619993968
$ws.Range("M4").Value = 29.427447
$ws.Range("N4").Value = 88.282341
$ws.Range("O4").Value = 0.1263788668529656
$ws.Range("P4").Value = 0.126503803280843
$ws.Range("Q4").Value = 298.5931096982085
$ws.Range("R4").Value = 1791.558658189251
$ws.Range("S4").Value = 0.01915908164708916
$ws.Range("T4").Value = 0.01394886115925935

$ws.Range("G5").Value = 10.1467555
$ws.Range("H5").Value = 20.293511
$ws.Range("I5").Value = 0.1516003594919049
$ws.Range("J5").Value = 0.1102643619993968
$ws.Range("M5").Value = 0.578241
$ws.Range("N5").Value = 1.156482
$ws.Range("O5").Value = 0.002483308944466901
$ws.Range("P5").Value = 0.001657175939929322
$ws.Range("Q5").Value = 5.8672700470755
$ws.Range("R5").Value = 23.469080188302
$ws.Range("S5").Value = 0.0003764705287106451
$ws.Range("T5").Value = 0.0001827274477370573

$ws.Range("G6").Value = 10.1467555
$ws.Range("H6").Value = 20.293511
$ws.Range("I6").Value = 0.1516003594919049
$ws.Range("J6").Value = 0.1102643619993968
$ws.Range("M6").Value = 194.164098
$ws.Range("N6").Value = 582.492294
$ws.Range("O6").Value = 0.8338555053303862
$ws.Range("P6").Value = 0.8346798435349938
$ws.Range("Q6").Value = 1970.135629284039
$ws.Range("R6").Value = 11820.81377570424
$ws.Range("S6").Value = 0.1264127943723906
$ws.Range("T6").Value = 0.0920354404211424

$ws.Range("G7").Value = 10.1467555
$ws.Range("H7").Value = 20.293511
$ws.Range("I7").Value = 0.1516003594919049
$ws.Range("J7").Value = 0.1102643619993968
$ws.Range("M7").Value = 8.221782333333335
$ws.Range("N7").Value = 24.665347
$ws.Range("O7").Value = 0.03530919738971574
$ws.Range("P7").Value = 0.03534410358173139
$ws.Range("Q7").Value = 83.42441511055286
$ws.Range("R7").Value = 500.5464906633171
$ws.Range("S7").Value = 0.005352887017651536
$ws.Range("T7").Value = 0.003897195031880205

$ws.Range("I8").Value = 0.1339722830802056
$ws.Range("J8").Value = 0.1461642475877201
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.111657
$ws.Range("N8").Value = 0.223314
$ws.Range("O8").Value = 0.0004795212148781231
$ws.Range("P8").Value = 0.0003199968420168896
$ws.Range("Q8").Value = 1.001216222825
$ws.Range("R8").Value = 6.007297336950001
$ws.Range("S8").Value = 0.00006424255194261599
$ws.Range("T8").Value = 0.00004677209764384521

$ws.Range("I9").Value = 0.1339722830802056
$ws.Range("J9").Value = 0.1461642475877201
$ws.Range("O9").Value = 0.001493600267587414
$ws.Range("P9").Value = 0.001495076820485504
$ws.Range("S9").Value = 0.0002001010378578918
$ws.Range("T9").Value = 0.0002185267785521046

$ws.Range("I10").Value = 0.1339722830802056
$ws.Range("J10").Value = 0.1461642475877201
$ws.Range("M10").Value = 29.427447
$ws.Range("N10").Value = 88.282341
$ws.Range("O10").Value = 0.1263788668529656
$ws.Range("P10").Value = 0.126503803280843
$ws.Range("Q10").Value = 263.872729275575
$ws.Range("R10").Value = 2374.854563480175
$ws.Range("S10").Value = 0.01693126532538111
$ws.Range("T10").Value = 0.01849033322352938

$ws.Range("I11").Value = 0.1339722830802056
$ws.Range("J11").Value = 0.1461642475877201
$ws.Range("M11").Value = 0.578241
$ws.Range("N11").Value = 1.156482
$ws.Range("O11").Value = 0.002483308944466901
$ws.Range("P11").Value = 0.001657175939929322
$ws.Range("Q11").Value = 5.185024404225
$ws.Range("R11").Value = 31.11014642535
$ws.Range("S11").Value = 0.0003326945688837261
$ws.Range("T11").Value = 0.0002422198743802422

$ws.Range("I12").Value = 0.1339722830802056
$ws.Range("J12").Value = 0.1461642475877201
$ws.Range("M12").Value = 194.164098
$ws.Range("N12").Value = 582.492294
$ws.Range("O12").Value = 0.8338555053303862
$ws.Range("P12").Value = 0.8346798435349938
$ws.Range("Q12").Value = 1741.04843232205
$ws.Range("R12").Value = 15669.43589089845
$ws.Range("S12").Value = 0.1117135258081104
$ws.Range("T12").Value = 0.1220003513069283

$ws.Range("I13").Value = 0.1339722830802056
$ws.Range("J13").Value = 0.1461642475877201
$ws.Range("M13").Value = 8.221782333333335
$ws.Range("N13").Value = 24.665347
$ws.Range("O13").Value = 0.03530919738971574
$ws.Range("P13").Value = 0.03534410358173139
$ws.Range("Q13").Value = 73.7238314899139
$ws.Range("R13").Value = 663.5144834092251
$ws.Range("S13").Value = 0.004730453788029852
$ws.Range("T13").Value = 0.005166044306686212

$ws.Range("G14").Value = 17.36323866666666
$ws.Range("H14").Value = 52.089716
$ws.Range("I14").Value = 0.2594201884346587
$ws.Range("J14").Value = 0.2830283681059314
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.111657
$ws.Range("N14").Value = 0.223314
$ws.Range("O14").Value = 0.0004795212148781231
$ws.Range("P14").Value = 0.0003199968420168896
$ws.Range("Q14").Value = 1.938727139804
$ws.Range("R14").Value = 11.632362838824
$ws.Range("S14").Value = 0.0001243974839220991
$ws.Range("T14").Value = 0.0000905681839950918

$ws.Range("G15").Value = 17.36323866666666
$ws.Range("H15").Value = 52.089716
$ws.Range("I15").Value = 0.2594201884346587
$ws.Range("J15").Value = 0.2830283681059314
$ws.Range("O15").Value = 0.001493600267587414
$ws.Range("P15").Value = 0.001495076820485504
$ws.Range("Q15").Value = 6.038697110671555
$ws.Range("R15").Value = 54.34827399604399
$ws.Range("S15").Value = 0.0003874700628635835
$ws.Range("T15").Value = 0.0004231491526950168

$ws.Range("G16").Value = 17.36323866666666
$ws.Range("H16").Value = 52.089716
$ws.Range("I16").Value = 0.2594201884346587
$ws.Range("J16").Value = 0.2830283681059314
$ws.Range("M16").Value = 29.427447
$ws.Range("N16").Value = 88.282341
$ws.Range("O16").Value = 0.1263788668529656
$ws.Range("P16").Value = 0.126503803280843
$ws.Range("Q16").Value = 510.955785611684
$ws.Range("R16").Value = 4598.602070505156
$ws.Range("S16").Value = 0.03278522945315497
$ws.Range("T16").Value = 0.03580416500177077

$ws.Range("G17").Value = 17.36323866666666
$ws.Range("H17").Value = 52.089716
$ws.Range("I17").Value = 0.2594201884346587
$ws.Range("J17").Value = 0.2830283681059314
$ws.Range("M17").Value = 0.578241
$ws.Range("N17").Value = 1.156482
$ws.Range("O17").Value = 0.002483308944466901
$ws.Range("P17").Value = 0.001657175939929322
$ws.Range("Q17").Value = 10.040136489852
$ws.Range("R17").Value = 60.240818939112
$ws.Range("S17").Value = 0.0006442204743150768
$ws.Range("T17").Value = 0.0004690278019426089

$ws.Range("G18").Value = 17.36323866666666
$ws.Range("H18").Value = 52.089716
$ws.Range("I18").Value = 0.2594201884346587
$ws.Range("J18").Value = 0.2830283681059314
$ws.Range("M18").Value = 194.164098
$ws.Range("N18").Value = 582.492294
$ws.Range("O18").Value = 0.8338555053303862
$ws.Range("P18").Value = 0.8346798435349938
$ws.Range("Q18").Value = 3371.317574072055
$ws.Range("R18").Value = 30341.8581666485
$ws.Range("S18").Value = 0.2163189523200863
$ws.Range("T18").Value = 0.2362380740066234

$ws.Range("G19").Value = 17.36323866666666
$ws.Range("H19").Value = 52.089716
$ws.Range("I19").Value = 0.2594201884346587
$ws.Range("J19").Value = 0.2830283681059314
$ws.Range("M19").Value = 8.221782333333335
$ws.Range("N19").Value = 24.665347
$ws.Range("O19").Value = 0.03530919738971574
$ws.Range("P19").Value = 0.03534410358173139
$ws.Range("Q19").Value = 142.7567689190502
$ws.Range("R19").Value = 1284.810920271452
$ws.Range("S19").Value = 0.009159918640316615
$ws.Range("T19").Value = 0.01000338395890444

$ws.Range("G20").Value = 6.601931
$ws.Range("H20").Value = 13.203862
$ws.Range("I20").Value = 0.09863794519743292
$ws.Range("J20").Value = 0.07174290438742112
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 0.111657
$ws.Range("N20").Value = 0.223314
$ws.Range("O20").Value = 0.0004795212148781231
$ws.Range("P20").Value = 0.0003199968420168896
$ws.Range("Q20").Value = 0.7371518096670001
$ws.Range("R20").Value = 2.948607238668
$ws.Range("S20").Value = 0.00004729898731415476
$ws.Range("T20").Value = 0.00002295750284109441

$ws.Range("G21").Value = 6.601931
$ws.Range("H21").Value = 13.203862
$ws.Range("I21").Value = 0.09863794519743292
$ws.Range("J21").Value = 0.07174290438742112
$ws.Range("O21").Value = 0.001493600267587414
$ws.Range("P21").Value = 0.001495076820485504
$ws.Range("Q21").Value = 2.296061375409667
$ws.Range("R21").Value = 13.776368252458
$ws.Range("S21").Value = 0.0001473256613411584
$ws.Range("T21").Value = 0.0001072611533839411

$ws.Range("G22").Value = 6.601931
$ws.Range("H22").Value = 13.203862
$ws.Range("I22").Value = 0.09863794519743292
$ws.Range("J22").Value = 0.07174290438742112
$ws.Range("M22").Value = 29.427447
$ws.Range("N22").Value = 88.282341
$ws.Range("O22").Value = 0.1263788668529656
$ws.Range("P22").Value = 0.126503803280843
$ws.Range("Q22").Value = 194.277974600157
$ws.Range("R22").Value = 1165.667847600942
$ws.Range("S22").Value = 0.01246575174275649
$ws.Range("T22").Value = 0.00907575026342265

$ws.Range("G23").Value = 6.601931
$ws.Range("H23").Value = 13.203862
$ws.Range("I23").Value = 0.09863794519743292
$ws.Range("J23").Value = 0.07174290438742112
$ws.Range("M23").Value = 0.578241
$ws.Range("N23").Value = 1.156482
$ws.Range("O23").Value = 0.002483308944466901
$ws.Range("P23").Value = 0.001657175939929322
$ws.Range("Q23").Value = 3.817507183371
$ws.Range("R23").Value = 15.270028733484
$ws.Range("S23").Value = 0.0002449484915726212
$ws.Range("T23").Value = 0.000118890615011484

$ws.Range("G24").Value = 6.601931
$ws.Range("H24").Value = 13.203862
$ws.Range("I24").Value = 0.09863794519743292
$ws.Range("J24").Value = 0.07174290438742112
$ws.Range("M24").Value = 194.164098
$ws.Range("N24").Value = 582.492294
$ws.Range("O24").Value = 0.8338555053303862
$ws.Range("P24").Value = 0.8346798435349938
$ws.Range("Q24").Value = 1281.857977673238
$ws.Range("R24").Value = 7691.147866039429
$ws.Range("S24").Value = 0.08224979363735638
$ws.Range("T24").Value = 0.05988235620883868

$ws.Range("G25").Value = 6.601931
$ws.Range("H25").Value = 13.203862
$ws.Range("I25").Value = 0.09863794519743292
$ws.Range("J25").Value = 0.07174290438742112
$ws.Range("M25").Value = 8.221782333333335
$ws.Range("N25").Value = 24.665347
$ws.Range("O25").Value = 0.03530919738971574
$ws.Range("P25").Value = 0.03534410358173139
$ws.Range("Q25").Value = 54.27963966168569
$ws.Range("R25").Value = 325.6778379701141
$ws.Range("S25").Value = 0.003482826677092123
$ws.Range("T25").Value = 0.002535688643923263

$ws.Range("G26").Value = 5.391932333333334
$ws.Range("H26").Value = 16.175797
$ws.Range("I26").Value = 0.080559631114533
$ws.Range("J26").Value = 0.0878908502346763
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0.111657
$ws.Range("N26").Value = 0.223314
$ws.Range("O26").Value = 0.0004795212148781231
$ws.Range("P26").Value = 0.0003199968420168896
$ws.Range("Q26").Value = 0.6020469885430001
$ws.Range("R26").Value = 3.612281931258001
$ws.Range("S26").Value = 0.00003863005218217431
$ws.Range("T26").Value = 0.00002812479451727582

$ws.Range("G27").Value = 5.391932333333334
$ws.Range("H27").Value = 16.175797
$ws.Range("I27").Value = 0.080559631114533
$ws.Range("J27").Value = 0.0878908502346763
$ws.Range("O27").Value = 0.001493600267587414
$ws.Range("P27").Value = 0.001495076820485504
$ws.Range("Q27").Value = 1.875240375791445
$ws.Range("R27").Value = 16.877163382123
$ws.Range("S27").Value = 0.0001203238865894098
$ws.Range("T27").Value = 0.0001314035729186275

$ws.Range("G28").Value = 5.391932333333334
$ws.Range("H28").Value = 16.175797
$ws.Range("I28").Value = 0.080559631114533
$ws.Range("J28").Value = 0.0878908502346763
$ws.Range("M28").Value = 29.427447
$ws.Range("N28").Value = 88.282341
$ws.Range("O28").Value = 0.1263788668529656
$ws.Range("P28").Value = 0.126503803280843
$ws.Range("Q28").Value = 158.670802966753
$ws.Range("R28").Value = 1428.037226700777
$ws.Range("S28").Value = 0.01018103489434759
$ws.Range("T28").Value = 0.01111852682827353

$ws.Range("G29").Value = 5.391932333333334
$ws.Range("H29").Value = 16.175797
$ws.Range("I29").Value = 0.080559631114533
$ws.Range("J29").Value = 0.0878908502346763
$ws.Range("M29").Value = 0.578241
$ws.Range("N29").Value = 1.156482
$ws.Range("O29").Value = 0.002483308944466901
$ws.Range("P29").Value = 0.001657175939929322
$ws.Range("Q29").Value = 3.117836344359
$ws.Range("R29").Value = 18.707018066154
$ws.Range("S29").Value = 0.0002000544525096738
$ws.Range("T29").Value = 0.0001456506023488369

$ws.Range("G30").Value = 5.391932333333334
$ws.Range("H30").Value = 16.175797
$ws.Range("I30").Value = 0.080559631114533
$ws.Range("J30").Value = 0.0878908502346763
$ws.Range("M30").Value = 194.164098
$ws.Range("N30").Value = 582.492294
$ws.Range("O30").Value = 0.8338555053303862
$ws.Range("P30").Value = 0.8346798435349938
$ws.Range("Q30").Value = 1046.919677978702
$ws.Range("R30").Value = 9422.27710180832
$ws.Range("S30").Value = 0.06717509191223842
$ws.Range("T30").Value = 0.07336072112203719

$ws.Range("G31").Value = 5.391932333333334
$ws.Range("H31").Value = 16.175797
$ws.Range("I31").Value = 0.080559631114533
$ws.Range("J31").Value = 0.0878908502346763
$ws.Range("M31").Value = 8.221782333333335
$ws.Range("N31").Value = 24.665347
$ws.Range("O31").Value = 0.03530919738971574
$ws.Range("P31").Value = 0.03534410358173139
$ws.Range("Q31").Value = 44.3312940007288
$ws.Range("R31").Value = 398.9816460065591
$ws.Range("S31").Value = 0.002844495916665731
$ws.Range("T31").Value = 0.003106423314580839

$ws.Range("G32").Value = 18.460197
$ws.Range("H32").Value = 55.380591
$ws.Range("I32").Value = 0.2758095926812649
$ws.Range("J32").Value = 0.3009092676848542
$ws.Range("K32").Value = 2
$ws.Range("L32").Value = 1
$ws.Range("M32").Value = 0.111657
$ws.Range("N32").Value = 0.223314
$ws.Range("O32").Value = 0.0004795212148781231
$ws.Range("P32").Value = 0.0003199968420168896
$ws.Range("Q32").Value = 2.061210216429
$ws.Range("R32").Value = 12.367261298574
$ws.Range("S32").Value = 0.0001322565509575604
$ws.Range("T32").Value = 0.00009629001539276824

$ws.Range("G33").Value = 18.460197
$ws.Range("H33").Value = 55.380591
$ws.Range("I33").Value = 0.2758095926812649
$ws.Range("J33").Value = 0.3009092676848542
$ws.Range("O33").Value = 0.001493600267587414
$ws.Range("P33").Value = 0.001495076820485504
$ws.Range("Q33").Value = 6.420204227241
$ws.Range("R33").Value = 57.781838045169
$ws.Range("S33").Value = 0.0004119492814319128
$ws.Range("T33").Value = 0.0004498824711848933

$ws.Range("G34").Value = 18.460197
$ws.Range("H34").Value = 55.380591
$ws.Range("I34").Value = 0.2758095926812649
$ws.Range("J34").Value = 0.3009092676848542
$ws.Range("M34").Value = 29.427447
$ws.Range("N34").Value = 88.282341
$ws.Range("O34").Value = 0.1263788668529656
$ws.Range("P34").Value = 0.126503803280843
$ws.Range("Q34").Value = 543.2364688270591
$ws.Range("R34").Value = 4889.128219443532
$ws.Range("S34").Value = 0.03485650379023624
$ws.Range("T34").Value = 0.03806616680458733

$ws.Range("G35").Value = 18.460197
$ws.Range("H35").Value = 55.380591
$ws.Range("I35").Value = 0.2758095926812649
$ws.Range("J35").Value = 0.3009092676848542
$ws.Range("M35").Value = 0.578241
$ws.Range("N35").Value = 1.156482
$ws.Range("O35").Value = 0.002483308944466901
$ws.Range("P35").Value = 0.001657175939929322
$ws.Range("Q35").Value = 10.674442773477
$ws.Range("R35").Value = 64.046656640862
$ws.Range("S35").Value = 0.0006849204284751578
$ws.Range("T35").Value = 0.000498659598509092

$ws.Range("G36").Value = 18.460197
$ws.Range("H36").Value = 55.380591
$ws.Range("I36").Value = 0.2758095926812649
$ws.Range("J36").Value = 0.3009092676848542
$ws.Range("M36").Value = 194.164098
$ws.Range("N36").Value = 582.492294
$ws.Range("O36").Value = 0.8338555053303862
$ws.Range("P36").Value = 0.8346798435349938
$ws.Range("Q36").Value = 3584.307499407306
$ws.Range("R36").Value = 32258.76749466576
$ws.Range("S36").Value = 0.2299853472802041
$ws.Range("T36").Value = 0.2511629004694237

$ws.Range("G37").Value = 18.460197
$ws.Range("H37").Value = 55.380591
$ws.Range("I37").Value = 0.2758095926812649
$ws.Range("J37").Value = 0.3009092676848542
$ws.Range("M37").Value = 8.221782333333335
$ws.Range("N37").Value = 24.665347
$ws.Range("O37").Value = 0.03530919738971574
$ws.Range("P37").Value = 0.03534410358173139
$ws.Range("Q37").Value = 151.7757215644531
$ws.Range("R37").Value = 1365.981494080077
$ws.Range("S37").Value = 0.00973861534995988
$ws.Range("T37").Value = 0.01063536832575642

